# Restore original multi-industry template text (undo "Product" find/replace
# that had been incorrectly applied over "AI"/"ML" terms), per commit 168d9c4.

$wb = $excel.ActiveWorkbook

# --- Sheet: Resource Overview ---
$ws1 = $wb.Worksheets.Item("Resource Overview")
$ws1.Range("A2").Value = "PRODUCT Resource Plan Staffing Plan Project"
$ws1.Range("B6").Value = "Enterprise AI/ML Implementation"
$ws1.Range("A18").Value = "Data Science/AI"
$ws1.Range("G18").Value = "ML, Python, Statistics"

# --- Sheet: Detailed Staffing Plan ---
$ws2 = $wb.Worksheets.Item("Detailed Staffing Plan")
$ws2.Range("A1").Value = "DETAILED STAFFING PLAN"
$ws2.Range("C9").Value = "Data Science/AI"
$ws2.Range("K9").Value = "ML, Deep Learning, Python"
$ws2.Range("P9").Value = "AI Lead"
$ws2.Range("C10").Value = "Data Science/AI"
$ws2.Range("K10").Value = "ML, Statistics, R/Python"
$ws2.Range("C11").Value = "Data Science/AI"
$ws2.Range("K11").Value = "ML, Python, Visualization"
$ws2.Range("B12").Value = "ML Engineer"
$ws2.Range("C12").Value = "Data Science/AI"
$ws2.Range("K12").Value = "MLOps, Python, Cloud"
$ws2.Range("C13").Value = "Data Science/AI"

# --- Sheet: Skills Matrix ---
$ws4 = $wb.Worksheets.Item("Skills Matrix")
$ws4.Range("D3").Value = "Machine Learning"

# --- Sheet: Cost Analysis ---
$ws5 = $wb.Worksheets.Item("Cost Analysis")
$ws5.Range("A6").Value = "Data Science/AI"

# --- Sheet: Resource Risk Assessment ---
$ws6 = $wb.Worksheets.Item("Resource Risk Assessment")
$ws6.Range("B5").Value = "Team lacks required ML expertise"
